$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)

# Reposition / resize the subtitle placeholder (was inheriting from layout,
# now gets an explicit xfrm matching the post-edit autofit layout).
$sh.Left = 82.80000305175781
$sh.Top = 362.5065612792969
$sh.Width = 621.3600463867188
$sh.Height = 84.24000549316406

# Turn on "shrink text on overflow" so PowerPoint stores <a:normAutofit/>
# on the body instead of leaving autosize unset.
$sh.TextFrame.AutoSize = 2

$tr = $sh.TextFrame.TextRange
$tr.Text = "Chan Choi"

# Blank paragraph, then the repo link on its own paragraph split across two
# runs (mirrors how the author typed the URL then the repo slug).
$tr.InsertAfter("`r`rgithub.com/chanchoi829/") | Out-Null
$tr2 = $sh.TextFrame.TextRange
$tr2.InsertAfter("battleship_cplusplus") | Out-Null
